# SAP_Test_Planner_Healthcheck.xlsx - "New File PP Regression6"
#
# Turn off the three VA_CO_0x rows (VA_CO_03 / VA_CO_01 / VA_CO_02 on rows
# 17-19) so they are skipped by the regression run, and tidy up the leftover
# shaded-fill formatting on the NC_OP_16..NC_OP_21 rows (20-25) in column C
# so they match the plain bordered style used elsewhere in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# Runmode column (B) - switch rows 17-19 from "Yes" to "No"
$ws.Range("B17:B19").Value = "No"

# Remove the stray fill/shading left on C20:C25 (TC_id column) so the cells
# go back to the plain no-fill bordered look used by the rest of the sheet.
$ws.Range("C20:C25").Interior.Pattern = -4142
